# MCS student tracking 2223.xlsx
#
# The author replaced the short placeholder labels in column A (rows 47-51)
# of the "feedback approval form" sheet with the fully written-out SMART
# goal sentences, and cleared the (now redundant) duplicate labels that used
# to sit in column B next to them. Finally the sheet selection was left on
# cell B51, scrolled down to around row 45 - which is where the edit was
# made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("feedback approval form")

# Row 47 - "Specific:" -> full sentence, clear duplicate in B47
$ws.Range("A47").Value = "Specific: Ik ben van plan voor het eind van deze sprint meerdere tests maken binnen de yml en dat dan over de code die we nu hebben te laten gaan"
$ws.Range("B47").ClearContents()

# Row 48 - "Measurable:" -> full sentence, clear duplicate in B48
$ws.Range("A48").Value = "Measurable: dit wordt dan measurable door de tests te zien in de pipeline"
$ws.Range("B48").ClearContents()

# Row 49 - "Achievable:" -> full sentence, clear duplicate in B49
$ws.Range("A49").Value = "Achievable: dit is haalbaar aangezien we genoeg code hebben, en er genoeg documentatie online staat"
$ws.Range("B49").ClearContents()

# Row 50 - "Realistic" -> full sentence, clear duplicate in B50
$ws.Range("A50").Value = "Realistic: dit is realistisch omdat het niet overdreven ingewikkelde code is en daarbij mogelijkheid is om het in de pipeline te laten testen"
$ws.Range("B50").ClearContents()

# Row 51 - "Time:" -> full sentence, clear duplicate in B51
$ws.Range("A51").Value = "Time: ik heb hier de tijd van de sprint van 3 weken genomen"
$ws.Range("B51").ClearContents()

# Leave the view/selection where the author ended up: scrolled to row 45,
# with B51 selected as the active cell.
$ws.Activate() | Out-Null
$ws.Range("B51").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
